$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.242.70"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "1.645.27"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'217.01"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "'19.96"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").Value = "'0.0792"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.873.93"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.30"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "1.656.99"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "'63.57"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "26.232.93"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'195.88"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "'10.04"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "'6.35"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "'0.125"
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "'0.0504"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "1.136.71"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "'0.553"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "'5.66"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "'100.22"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").Value = "1.782.86"
$ws.Range("D46").Value = "'56.26"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("D49").Value = "'0.418"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("D51").Value = "'0.0975"
$ws.Range("E51").Value = "  +2.30%  "
